$d = $word.ActiveDocument

# 1) The word "Version" is currently split across two runs ("Versi" + "on").
#    Re-assert the same text over the whole word so Word collapses it back
#    into a single run, matching the reverted (pre-typo-fix) formatting.
$rngVersion = $d.Content
$rngVersion.Find.Execute("Version", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Version", 2) | Out-Null

# 2) Change the version number from " 2" to " 1."
$rngNum = $d.Content
$rngNum.Find.Execute(" 2", $true, $false, $false, $false, $false, `
                      $true, 1, $false, " 1.", 2) | Out-Null

# 3) The trailing sentence-ending "." used to live in its own run after the
#    bookmark; that period is now already part of the " 1." run above, so
#    drop the now-duplicate trailing "." character (the one right before
#    the paragraph mark, after the bookmark).
$p = $d.Paragraphs(1)
$tail = $d.Range($p.Range.End - 2, $p.Range.End - 1)
if ($tail.Text -eq ".") {
    $tail.Text = ""
}
